# Rename the "_old"/"_new" header-suffix scheme to the format-version
# specific suffixes "_FV2404" (old/before) and "_FV2410" (new/after),
# turn the header range into a real Excel Table ("Table1"), and freeze
# the header row (row 1) in the view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2404", "Segmentgruppe_FV2404", "Segment_FV2404", "Datenelement_FV2404", "Segment ID_FV2404",
    "Code_FV2404", "Qualifier_FV2404", "Beschreibung_FV2404", "Bedingungsausdruck_FV2404", "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410", "Segmentgruppe_FV2410", "Segment_FV2410", "Datenelement_FV2410", "Segment ID_FV2410",
    "Code_FV2410", "Qualifier_FV2410", "Beschreibung_FV2410", "Bedingungsausdruck_FV2410", "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Turn A1:U53 into a proper Excel Table named "Table1" (adds
# xl/tables/table1.xml, the worksheet <tableParts> reference and the
# relationship part).
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U53"), $null, 1)
$lo.Name = "Table1"

# Freeze the header row so row 1 stays visible while scrolling.
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
